# Apply updated odds/statistics values for rows 4, 5 and 6 in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4
$ws.Range("G4").Value = 2.1
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 3.75
$ws.Range("Z4").Value = 19
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 21
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 19
$ws.Range("AN4").Value = 4.33
$ws.Range("AO4").Value = 11
$ws.Range("AQ4").Value = 34
$ws.Range("AS4").Value = 101
$ws.Range("AW4").Value = 5.5
$ws.Range("AY4").Value = 23
$ws.Range("BD4").Value = 151

# Row 5
$ws.Range("G5").Value = 1.53
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("AS5").Value = 201

# Row 6
$ws.Range("G6").Value = 3.3
$ws.Range("Z6").Value = 34
